# Updates the "cryptos" price/volume table with the latest scraped values.
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h). Rows 2-51 hold the data.
#
# Some of the new Price values (column D) are plain decimals such as
# "522.10" / "4.00" / "36.01" that Excel would otherwise auto-convert to a
# number (dropping the trailing zero / significant digits), whereas the
# source data stores them as plain text. Set-TextValue forces the cell to
# text (NumberFormat "@") before assigning the value, then resets the
# style back to "Normal" so no stray style index is left on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $cellRef, $val)
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "58.152.03"
$ws.Range("E2").Value = "  -0.37%  "
$ws.Range("D3").Value = "2.593.12"
$ws.Range("E3").Value = "  -0.86%  "
$ws.Range("E4").Value = "  +0.12%  "
Set-TextValue $ws "D5" "522.10"
$ws.Range("E5").Value = "  +0.36%  "
Set-TextValue $ws "D6" "143.27"
$ws.Range("E6").Value = "  +0.88%  "
$ws.Range("E7").Value = "  -0.11%  "
Set-TextValue $ws "D8" "0.568"
$ws.Range("E8").Value = "  +0.61%  "
$ws.Range("D9").Value = "2.612.37"
$ws.Range("E9").Value = "  -0.30%  "
Set-TextValue $ws "D10" "6.46"
$ws.Range("E10").Value = "  -1.31%  "
Set-TextValue $ws "D12" "0.342"
$ws.Range("E12").Value = "  +1.76%  "
$ws.Range("E13").Value = "  +0.01%  "
$ws.Range("D14").Value = "3.053.06"
$ws.Range("E14").Value = "  -0.61%  "
$ws.Range("D15").Value = "58.111.54"
$ws.Range("E15").Value = "  -0.36%  "
Set-TextValue $ws "D16" "20.36"
$ws.Range("E16").Value = "  -2.60%  "
$ws.Range("E17").Value = "  -1.15%  "
$ws.Range("D18").Value = "2.562.55"
$ws.Range("E18").Value = "  -2.41%  "
Set-TextValue $ws "D19" "339.42"
$ws.Range("E19").Value = "  +0.75%  "
Set-TextValue $ws "D20" "4.34"
$ws.Range("E20").Value = "  -1.08%  "
$ws.Range("E21").Value = "  -1.44%  "
Set-TextValue $ws "D22" "6.42"
$ws.Range("E22").Value = "  +2.25%  "
Set-TextValue $ws "D23" "0.997"
$ws.Range("E23").Value = "  -0.16%  "
Set-TextValue $ws "D24" "65.33"
$ws.Range("E24").Value = "  +0.80%  "
$ws.Range("E25").Value = "  +1.31%  "
$ws.Range("E26").Value = "  -2.20%  "
$ws.Range("D27").Value = "2.717.86"
$ws.Range("E27").Value = "  -0.56%  "
$ws.Range("E28").Value = "  -0.01%  "
Set-TextValue $ws "D29" "7.03"
$ws.Range("E29").Value = "  -1.10%  "
$ws.Range("D30").Value = "0.0₃0745"
$ws.Range("E31").Value = "  -0.05%  "
Set-TextValue $ws "D32" "6.14"
$ws.Range("E32").Value = "  -5.90%  "
$ws.Range("E33").Value = "  -0.32%  "
Set-TextValue $ws "D34" "18.77"
$ws.Range("E34").Value = "  +0.13%  "
Set-TextValue $ws "D35" "149.65"
$ws.Range("E35").Value = "  -0.01%  "
Set-TextValue $ws "D36" "4.00"
$ws.Range("E36").Value = "  -1.81%  "
$ws.Range("E37").Value = "  -3.88%  "
Set-TextValue $ws "D38" "0.864"
$ws.Range("E38").Value = "  -2.10%  "
Set-TextValue $ws "D39" "0.863"
$ws.Range("E39").Value = "  +1.67%  "
$ws.Range("E40").Value = "  +2.73%  "
Set-TextValue $ws "D41" "36.01"
$ws.Range("E41").Value = "  -0.68%  "
Set-TextValue $ws "D42" "3.54"
$ws.Range("E42").Value = "  -2.21%  "
$ws.Range("E43").Value = "  -0.29%  "
Set-TextValue $ws "D44" "0.608"
$ws.Range("E44").Value = "  +0.78%  "
$ws.Range("E45").Value = "  +1.10%  "
$ws.Range("E46").Value = "  -0.92%  "
$ws.Range("E47").Value = "  +0.19%  "
Set-TextValue $ws "D48" "18.74"
$ws.Range("E48").Value = "  -1.62%  "
Set-TextValue $ws "D49" "0.0522"
$ws.Range("E49").Value = "  -1.53%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "1.967.69"
$ws.Range("E50").Value = "  -2.82%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws "D51" "18.72"
$ws.Range("E51").Value = "  +2.97%  "
